$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E contain text-formatted values (prices / percentages).
# Force text format first so Excel does not auto-convert numeric-looking
# strings (e.g. "441.30") into real numbers when we assign them below.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "52.765.69"
$ws.Range("E2").Value = "  -12.58%  "
$ws.Range("D3").Value = "2.329.70"
$ws.Range("E3").Value = "  -19.37%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "441.30"
$ws.Range("E5").Value = "  -15.92%  "
$ws.Range("D6").Value = "122.61"
$ws.Range("E6").Value = "  -12.98%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "0.469"
$ws.Range("E8").Value = "  -14.24%  "
$ws.Range("D9").Value = "2.337.70"
$ws.Range("E9").Value = "  -19.22%  "
$ws.Range("D10").Value = "5.26"
$ws.Range("E10").Value = "  -11.83%  "
$ws.Range("D11").Value = "0.0883"
$ws.Range("E11").Value = "  -17.20%  "
$ws.Range("D12").Value = "0.305"
$ws.Range("E12").Value = "  -14.56%  "
$ws.Range("D13").Value = "0.120"
$ws.Range("E13").Value = "  -5.91%  "
$ws.Range("D14").Value = "52.753.58"
$ws.Range("E14").Value = "  -12.67%  "
$ws.Range("D15").Value = "19.00"
$ws.Range("E15").Value = "  -15.74%  "
$ws.Range("E16").Value = "  -14.97%  "
$ws.Range("D17").Value = "2.338.49"
$ws.Range("E17").Value = "  -18.96%  "
$ws.Range("D18").Value = "3.96"
$ws.Range("E18").Value = "  -20.34%  "
$ws.Range("D19").Value = "300.53"
$ws.Range("E19").Value = "  -14.90%  "
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").Value = "  -22.42%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "5.58"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  -21.39%  "
$ws.Range("D24").Value = "53.95"
$ws.Range("E24").Value = "  -16.16%  "
$ws.Range("D25").Value = "0.150"
$ws.Range("E25").Value = "  -15.72%  "
$ws.Range("D26").Value = "0.366"
$ws.Range("E26").Value = "  -18.88%  "
$ws.Range("D27").Value = "6.92"
$ws.Range("E27").Value = "  -11.30%  "
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "0.0₃0679"
$ws.Range("E29").Value = "  -17.98%  "
$ws.Range("D30").Value = "142.40"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("D31").Value = "17.07"
$ws.Range("E31").Value = "  -12.67%  "
$ws.Range("E32").Value = "  -19.44%  "
$ws.Range("D33").Value = "4.79"
$ws.Range("E33").Value = "  -13.54%  "
$ws.Range("E34").Value = "  -15.72%  "
$ws.Range("D35").Value = "3.46"
$ws.Range("E35").Value = "  -19.32%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -16.01%  "
$ws.Range("D37").Value = "0.995"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "31.93"
$ws.Range("E38").Value = "  -14.84%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "10.15"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0506"
$ws.Range("E40").Value = "  -12.62%  "
$ws.Range("D41").Value = "3.15"
$ws.Range("E41").Value = "  -14.77%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.552"
$ws.Range("E42").Value = "  -14.36%  "
$ws.Range("D43").Value = "1.925.56"
$ws.Range("E43").Value = "  -15.56%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.21"
$ws.Range("E44").Value = "  -17.32%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.25"
$ws.Range("E45").Value = "  -12.98%  "
$ws.Range("D46").Value = "0.0827"
$ws.Range("E46").Value = "  -9.74%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0206"
$ws.Range("E47").Value = "  -12.78%  "
$ws.Range("D48").Value = "15.65"
$ws.Range("E48").Value = "  -22.41%  "
$ws.Range("D49").Value = "4.57"
$ws.Range("E49").Value = "  -5.29%  "
$ws.Range("D50").Value = "4.48"
$ws.Range("E50").Value = "  -12.54%  "
$ws.Range("D51").Value = "15.11"
$ws.Range("E51").Value = "  -16.53%  "

# Restore the original (default) cell style now that the text values are set.
$ws.Range("D2:E51").Style = "Normal"
